$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $style = $range.Style
    $range.Value = "'" + $value
    $range.Style = $style
}

Set-TextValue $ws.Range("D2") "44.280.87"
Set-TextValue $ws.Range("E2") "  +4.86%  "

Set-TextValue $ws.Range("D3") "2.261.19"
Set-TextValue $ws.Range("E3") "  +2.07%  "

Set-TextValue $ws.Range("E4") "  -0.03%  "

Set-TextValue $ws.Range("D5") "230.06"
Set-TextValue $ws.Range("E5") "  -0.31%  "

Set-TextValue $ws.Range("D6") "0.632"
Set-TextValue $ws.Range("E6") "  +2.38%  "

Set-TextValue $ws.Range("D7") "63.43"
Set-TextValue $ws.Range("E7") "  +4.37%  "

Set-TextValue $ws.Range("E8") "  +0.05%  "

Set-TextValue $ws.Range("D9") "0.444"
Set-TextValue $ws.Range("E9") "  +10.52%  "

Set-TextValue $ws.Range("E10") "  +14.47%  "

Set-TextValue $ws.Range("D11") "57.09"
Set-TextValue $ws.Range("E11") "  -0.55%  "

Set-TextValue $ws.Range("D12") "25.97"
Set-TextValue $ws.Range("E12") "  +17.41%  "

Set-TextValue $ws.Range("D13") "0.105"
Set-TextValue $ws.Range("E13") "  +1.98%  "

Set-TextValue $ws.Range("D14") "2.598.09"
Set-TextValue $ws.Range("E14") "  +2.18%  "

Set-TextValue $ws.Range("D15") "15.64"
Set-TextValue $ws.Range("E15") "  +1.29%  "

Set-TextValue $ws.Range("D16") "6.16"
Set-TextValue $ws.Range("E16") "  +10.61%  "

Set-TextValue $ws.Range("D17") "0.834"
Set-TextValue $ws.Range("E17") "  +4.85%  "

Set-TextValue $ws.Range("D18") "2.280.58"
Set-TextValue $ws.Range("E18") "  +2.48%  "

Set-TextValue $ws.Range("D19") "43.970.98"
Set-TextValue $ws.Range("E19") "  +4.52%  "

Set-TextValue $ws.Range("E20") "  +8.91%  "

Set-TextValue $ws.Range("D21") "73.24"
Set-TextValue $ws.Range("E21") "  +1.74%  "

Set-TextValue $ws.Range("E22") "  -2.75%  "

Set-TextValue $ws.Range("D23") "251.07"
Set-TextValue $ws.Range("E23") "  +3.18%  "

Set-TextValue $ws.Range("E24") "  +0.05%  "

Set-TextValue $ws.Range("D25") "2.42"
Set-TextValue $ws.Range("E25") "  +0.39%  "

Set-TextValue $ws.Range("E26") "  -1.79%  "

$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D27") "9.99"
Set-TextValue $ws.Range("E27") "  +4.20%  "

$ws.Range("B28").Value = "WEMIXToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D28") "3.27"
Set-TextValue $ws.Range("E28") "  +22.70%  "

Set-TextValue $ws.Range("D29") "171.94"
Set-TextValue $ws.Range("E29") "  +1.77%  "

Set-TextValue $ws.Range("D30") "20.72"
Set-TextValue $ws.Range("E30") "  +1.94%  "

Set-TextValue $ws.Range("E31") "  -2.90%  "

Set-TextValue $ws.Range("E32") "  -5.68%  "

Set-TextValue $ws.Range("E33") "  +2.38%  "

Set-TextValue $ws.Range("D34") "0.0679"
Set-TextValue $ws.Range("E34") "  +4.88%  "

Set-TextValue $ws.Range("D35") "4.72"
Set-TextValue $ws.Range("E35") "  +2.54%  "

Set-TextValue $ws.Range("D36") "4.84"
Set-TextValue $ws.Range("E36") "  -2.41%  "

Set-TextValue $ws.Range("D37") "3.81"
Set-TextValue $ws.Range("E37") "  +7.29%  "

Set-TextValue $ws.Range("D38") "6.64"
Set-TextValue $ws.Range("E38") "  +5.55%  "

Set-TextValue $ws.Range("D39") "2.29"
Set-TextValue $ws.Range("E39") "  -1.58%  "

Set-TextValue $ws.Range("D40") "0.0258"
Set-TextValue $ws.Range("E40") "  +3.90%  "

Set-TextValue $ws.Range("E41") "  -0.08%  "

$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D42") "17.39"
Set-TextValue $ws.Range("E42") "  +8.53%  "

$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue $ws.Range("D43") "8.27"
Set-TextValue $ws.Range("E43") "  -3.04%  "

Set-TextValue $ws.Range("D44") "0.0963"
Set-TextValue $ws.Range("E44") "  +0.97%  "

Set-TextValue $ws.Range("D45") "97.30"
Set-TextValue $ws.Range("E45") "  +0.43%  "

Set-TextValue $ws.Range("D46") "1.18"
Set-TextValue $ws.Range("E46") "  -0.99%  "

$ws.Range("B47").Value = "TerraClassic"
$ws.Range("C47").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
Set-TextValue $ws.Range("D47") "0.000210"
Set-TextValue $ws.Range("E47") "  -7.42%  "

$ws.Range("B48").Value = "FTXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D48") "4.34"
Set-TextValue $ws.Range("E48") "  +0.37%  "

Set-TextValue $ws.Range("D49") "1.437.47"
Set-TextValue $ws.Range("E49") "  -1.24%  "

Set-TextValue $ws.Range("D50") "2.27"

Set-TextValue $ws.Range("D51") "2.75"
Set-TextValue $ws.Range("E51") "  +1.49%  "
